# Auto-generated edit script: refresh market-data snapshot values
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 56.625
$ws.Range("I6").Value = 59.142857
$ws.Range("K6").Value = 177.428571
$ws.Range("M6").Value = -65.42857100000001
$ws.Range("H19").Value = 1401
$ws.Range("I19").Value = 1459.1428
$ws.Range("K19").Value = 1459.1428
$ws.Range("M19").Value = -1284.1428
$ws.Range("H32").Value = 7000.5
$ws.Range("J32").Value = 7000.6
$ws.Range("L32").Value = 7000.6
$ws.Range("N32").Value = -7652.6
$ws.Range("H37").Value = 600
$ws.Range("I37").Value = 600
$ws.Range("K37").Value = 1800
$ws.Range("M37").Value = -1674
$ws.Range("H80").Value = 4683.8887
$ws.Range("I80").Value = 5184
$ws.Range("J80").Value = 4433.8335
$ws.Range("K80").Value = 15552
$ws.Range("L80").Value = 13301.5005
$ws.Range("M80").Value = -14554
$ws.Range("N80").Value = -15297.5005
$ws.Range("H83").Value = 4683.8887
$ws.Range("I83").Value = 5184
$ws.Range("J83").Value = 4433.8335
$ws.Range("K83").Value = 46656
$ws.Range("L83").Value = 39904.5015
$ws.Range("M83").Value = -41664
$ws.Range("N83").Value = -49888.5015
$ws.Range("H106").Value = 4998.3335
$ws.Range("I106").Value = 4998.3335
$ws.Range("K106").Value = 4998.3335
$ws.Range("M106").Value = -4367.3335
$ws.Range("H107").Value = 449.6
$ws.Range("I107").Value = 537
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 537
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 1383
$ws.Range("N107").Value = -3940
$ws.Range("H137").Value = 1273.8
$ws.Range("I137").Value = 1217.25
$ws.Range("K137").Value = 3651.75
$ws.Range("M137").Value = -1101.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4830.8
$ws.Range("I45").Value = 5399.75
$ws.Range("K45").Value = 5399.75
$ws.Range("M45").Value = -5022.75
$ws.Range("H101").Value = 300200.34
$ws.Range("J101").Value = 300200.34
$ws.Range("L101").Value = 300200.34
$ws.Range("N101").Value = -306690.34
$ws.Range("H102").Value = 999
$ws.Range("I102").Value = 999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 623
$ws.Range("N102").Value = $null
$ws.Range("H132").Value = 2507
$ws.Range("I132").Value = 2314.1667
$ws.Range("J132").Value = 3374.75
$ws.Range("K132").Value = 6942.500100000001
$ws.Range("L132").Value = 10124.25
$ws.Range("M132").Value = -4412.500100000001
$ws.Range("N132").Value = -15184.25
$ws.Range("H135").Value = 278498.75
$ws.Range("J135").Value = 278498.75
$ws.Range("L135").Value = 278498.75
$ws.Range("N135").Value = -288638.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 10718.667
$ws.Range("J103").Value = 10718.667
$ws.Range("L103").Value = 10718.667
$ws.Range("N103").Value = -13062.667
$ws.Range("H105").Value = 3096.6667
$ws.Range("I105").Value = 3357
$ws.Range("K105").Value = 3357
$ws.Range("M105").Value = -1610

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1244.5
$ws.Range("I22").Value = 989
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 989
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -639
$ws.Range("N22").Value = -2200
$ws.Range("H25").Value = 5200
$ws.Range("I25").Value = 5200
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 5200
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -5026
$ws.Range("N25").Value = $null
$ws.Range("H59").Value = 113714.29
$ws.Range("I59").Value = 95000
$ws.Range("J59").Value = 121200
$ws.Range("K59").Value = 95000
$ws.Range("L59").Value = 121200
$ws.Range("M59").Value = -93855
$ws.Range("N59").Value = -123490
$ws.Range("H86").Value = 9828
$ws.Range("I86").Value = 9759.6
$ws.Range("K86").Value = 9759.6
$ws.Range("M86").Value = -8636.6
$ws.Range("H89").Value = 9828
$ws.Range("I89").Value = 9759.6
$ws.Range("K89").Value = 48798
$ws.Range("M89").Value = -43182
$ws.Range("H92").Value = 48601
$ws.Range("J92").Value = 48601
$ws.Range("L92").Value = 48601
$ws.Range("N92").Value = -53593

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79966
$ws.Range("J37").Value = 79966
$ws.Range("L37").Value = 239898
$ws.Range("N37").Value = -240122
$ws.Range("H80").Value = 1419
$ws.Range("I80").Value = 1681
$ws.Range("J80").Value = 1244.3334
$ws.Range("K80").Value = 5043
$ws.Range("L80").Value = 3733.0002
$ws.Range("M80").Value = -4107
$ws.Range("N80").Value = -5605.0002
$ws.Range("H83").Value = 1419
$ws.Range("I83").Value = 1681
$ws.Range("J83").Value = 1244.3334
$ws.Range("K83").Value = 15129
$ws.Range("L83").Value = 11199.0006
$ws.Range("M83").Value = -10449
$ws.Range("N83").Value = -20559.0006
$ws.Range("H113").Value = 2383.3333
$ws.Range("J113").Value = 2383.3333
$ws.Range("L113").Value = 7149.999899999999
$ws.Range("N113").Value = -11489.9999
$ws.Range("H114").Value = 1939.0714
$ws.Range("I114").Value = 2121.5715
$ws.Range("J114").Value = 1756.5714
$ws.Range("K114").Value = 6364.7145
$ws.Range("L114").Value = 5269.7142
$ws.Range("M114").Value = -3110.7145
$ws.Range("N114").Value = -11777.7142
$ws.Range("H117").Value = 11646.777
$ws.Range("J117").Value = 25673.75
$ws.Range("L117").Value = 77021.25
$ws.Range("N117").Value = -83905.25
$ws.Range("H129").Value = 2651.6365
$ws.Range("J129").Value = 3058.5
$ws.Range("L129").Value = 9175.5
$ws.Range("N129").Value = -19175.5
$ws.Range("H130").Value = 1911.1666
$ws.Range("J130").Value = 1885.6666
$ws.Range("L130").Value = 5656.9998
$ws.Range("N130").Value = -15696.9998
$ws.Range("H132").Value = 2254.4285
$ws.Range("I132").Value = 592.3333
$ws.Range("K132").Value = 5330.9997
$ws.Range("M132").Value = -2800.9997
$ws.Range("H137").Value = 1819.9166
$ws.Range("I137").Value = 1434.1
$ws.Range("J137").Value = 3749
$ws.Range("K137").Value = 4302.299999999999
$ws.Range("L137").Value = 11247
$ws.Range("M137").Value = 797.7000000000007
$ws.Range("N137").Value = -21447
$ws.Range("H139").Value = 2497
$ws.Range("I139").Value = 2497
$ws.Range("K139").Value = 7491
$ws.Range("M139").Value = -2351

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 844
$ws.Range("I97").Value = 874.8
$ws.Range("J97").Value = 690
$ws.Range("K97").Value = 874.8
$ws.Range("L97").Value = 690
$ws.Range("M97").Value = -378.8
$ws.Range("N97").Value = -1682
$ws.Range("H102").Value = 1064.7333
$ws.Range("I102").Value = 1090.4615
$ws.Range("J102").Value = 897.5
$ws.Range("K102").Value = 1090.4615
$ws.Range("L102").Value = 897.5
$ws.Range("M102").Value = 531.5385000000001
$ws.Range("N102").Value = -4141.5
$ws.Range("H126").Value = 4497.5
$ws.Range("I126").Value = 4195
$ws.Range("K126").Value = 12585
$ws.Range("M126").Value = -10115

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6836.407
$ws.Range("I7").Value = 2499.625
$ws.Range("J7").Value = 8662.421
$ws.Range("K7").Value = 2499.625
$ws.Range("L7").Value = 8662.421
$ws.Range("M7").Value = -2387.625
$ws.Range("N7").Value = -8886.421
$ws.Range("H16").Value = 2164.5557
$ws.Range("I16").Value = 958.4
$ws.Range("K16").Value = 958.4
$ws.Range("M16").Value = -788.4
$ws.Range("H40").Value = 3005.3845
$ws.Range("I40").Value = 1710.1428
$ws.Range("K40").Value = 1710.1428
$ws.Range("M40").Value = -1574.1428
$ws.Range("H68").Value = 2997.5
$ws.Range("J68").Value = 2997.5
$ws.Range("L68").Value = 2997.5
$ws.Range("N68").Value = -4495.5
$ws.Range("H71").Value = 2997.5
$ws.Range("J71").Value = 2997.5
$ws.Range("L71").Value = 14987.5
$ws.Range("N71").Value = -22475.5
$ws.Range("H122").Value = 7742.3335
$ws.Range("I122").Value = 8113.9287
$ws.Range("J122").Value = 6999.143
$ws.Range("K122").Value = 24341.7861
$ws.Range("L122").Value = 20997.429
$ws.Range("M122").Value = -21891.7861
$ws.Range("N122").Value = -25897.429
$ws.Range("H126").Value = 6836.407
$ws.Range("I126").Value = 2499.625
$ws.Range("J126").Value = 8662.421
$ws.Range("K126").Value = 7498.875
$ws.Range("L126").Value = 25987.263
$ws.Range("M126").Value = -5028.875
$ws.Range("N126").Value = -30927.263
$ws.Range("H132").Value = 2920.762
$ws.Range("I132").Value = 2795.6667
$ws.Range("J132").Value = 3671.3333
$ws.Range("K132").Value = 8387.000100000001
$ws.Range("L132").Value = 11013.9999
$ws.Range("M132").Value = -5857.000100000001
$ws.Range("N132").Value = -16073.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18618
$ws.Range("I62").Value = 24835
$ws.Range("J62").Value = 14732.375
$ws.Range("K62").Value = 24835
$ws.Range("L62").Value = 14732.375
$ws.Range("M62").Value = -24211
$ws.Range("N62").Value = -15980.375
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H65").Value = 18618
$ws.Range("I65").Value = 24835
$ws.Range("J65").Value = 14732.375
$ws.Range("K65").Value = 124175
$ws.Range("L65").Value = 73661.875
$ws.Range("M65").Value = -121055
$ws.Range("N65").Value = -79901.875
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H122").Value = 1129.5
$ws.Range("I122").Value = 923.5
$ws.Range("K122").Value = 2770.5
$ws.Range("M122").Value = -320.5
$ws.Range("H132").Value = 9200.200000000001
$ws.Range("I132").Value = 10250.25
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 30750.75
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -28220.75
$ws.Range("N132").Value = -20060
